# Update demo data to be more realistic and complete
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update folder_id values on existing rows (2-5) ---
$ws.Range("B2").Value = "13-protec-social"
$ws.Range("B3").Value = "00-base"
$ws.Range("B4").Value = "00-base"
$ws.Range("B5").Value = "00-base"

# --- Add new rows 6-9 ---
$ws.Range("A6").Value = "canton_sigle"
$ws.Range("B6").Value = "00-base"
$ws.Range("C6").Value = "Sigle des cantons"
$ws.Range("D6").Value = "Liste des 26 cantons suisses et leur sigle"

$ws.Range("A7").Value = "langue_sigle"
$ws.Range("B7").Value = "00-base"
$ws.Range("C7").Value = "Sigle des langues"
$ws.Range("D7").Value = "Liste des 3 principales langues suisses et leur sigle"

$ws.Range("A8").Value = "oui_non"
$ws.Range("B8").Value = "00-base"
$ws.Range("C8").Value = "Oui ou non"
$ws.Range("D8").Value = "Oui ou non codé en 0 ou 1"

$ws.Range("A9").Value = "vide"
$ws.Range("B9").Value = "00-base"
$ws.Range("C9").Value = "vide / manquant"
$ws.Range("D9").Value = "Valeur vide ou manquante"

# --- Apply left/center alignment style to A6:A8 (matches new cellXfs style) ---
$a6 = $ws.Range("A6")
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108
$a6.Copy() | Out-Null
$ws.Range("A7:A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Resize the table to include new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D9"))

# --- Resize columns C and D to fit new (wider) content ---
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 39.1

# --- Update selection to reflect new active cell ---
$ws.Range("D10").Select() | Out-Null
